$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213, shifting existing rows 213..296 down to 214..297
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new data record
$ws.Range("A213").Value = 4
$ws.Range("B213").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value = "Los Lagos"
$ws.Range("D213").Value = 44468
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = 100112004
$ws.Range("G213").Value = "Cebolla"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "1a (guarda)"
$ws.Range("J213").Value = 200
$ws.Range("K213").Value = 7000
$ws.Range("L213").Value = 7000
$ws.Range("M213").Value = 7000
$ws.Range("N213").Value = "`$/malla 16 kilos"
$ws.Range("O213").Value = "Región de O'Higgins"
$ws.Range("P213").Value = 438
$ws.Range("Q213").Value = 16
$ws.Range("R213").Value = "Hortaliza"
